$d = $word.ActiveDocument

$d.Content.Find.Execute("69-20=49", $true, $false, $false, $false, $false, $true, 1, $false, "58-19=39", 2) | Out-Null
$d.Content.Find.Execute("8+48=56", $true, $false, $false, $false, $false, $true, 1, $false, "69-24=45", 2) | Out-Null
$d.Content.Find.Execute("67-56=11", $true, $false, $false, $false, $false, $true, 1, $false, "31+1=32", 2) | Out-Null
$d.Content.Find.Execute("90-12=78", $true, $false, $false, $false, $false, $true, 1, $false, "24+40=64", 2) | Out-Null
$d.Content.Find.Execute("90-82=8", $true, $false, $false, $false, $false, $true, 1, $false, "94-87=7", 2) | Out-Null
$d.Content.Find.Execute("66-45=21", $true, $false, $false, $false, $false, $true, 1, $false, "50+35=85", 2) | Out-Null
$d.Content.Find.Execute("10+26=36", $true, $false, $false, $false, $false, $true, 1, $false, "12+28=40", 2) | Out-Null
$d.Content.Find.Execute("71-59=12", $true, $false, $false, $false, $false, $true, 1, $false, "9+3=12", 2) | Out-Null
$d.Content.Find.Execute("56-45=11", $true, $false, $false, $false, $false, $true, 1, $false, "85-64=21", 2) | Out-Null
$d.Content.Find.Execute("30+39=69", $true, $false, $false, $false, $false, $true, 1, $false, "21+12=33", 2) | Out-Null
$d.Content.Find.Execute("77+14=91", $true, $false, $false, $false, $false, $true, 1, $false, "40+50=90", 2) | Out-Null
$d.Content.Find.Execute("97-88=9", $true, $false, $false, $false, $false, $true, 1, $false, "74-36=38", 2) | Out-Null
$d.Content.Find.Execute("19+63=82", $true, $false, $false, $false, $false, $true, 1, $false, "5+19=24", 2) | Out-Null
$d.Content.Find.Execute("34+25=59", $true, $false, $false, $false, $false, $true, 1, $false, "77+3=80", 2) | Out-Null
$d.Content.Find.Execute("46+18=64", $true, $false, $false, $false, $false, $true, 1, $false, "15+29=44", 2) | Out-Null
$d.Content.Find.Execute("42+0=42", $true, $false, $false, $false, $false, $true, 1, $false, "27+13=40", 2) | Out-Null
$d.Content.Find.Execute("91-66=25", $true, $false, $false, $false, $false, $true, 1, $false, "74-38=36", 2) | Out-Null
$d.Content.Find.Execute("11+53=64", $true, $false, $false, $false, $false, $true, 1, $false, "47+48=95", 2) | Out-Null
$d.Content.Find.Execute("86-7=79", $true, $false, $false, $false, $false, $true, 1, $false, "6+56=62", 2) | Out-Null
$d.Content.Find.Execute("2+53=55", $true, $false, $false, $false, $false, $true, 1, $false, "23+37=60", 2) | Out-Null
$d.Content.Find.Execute("48-37=11", $true, $false, $false, $false, $false, $true, 1, $false, "0+67=67", 2) | Out-Null
$d.Content.Find.Execute("70-15=55", $true, $false, $false, $false, $false, $true, 1, $false, "35+9=44", 2) | Out-Null
$d.Content.Find.Execute("43-29=14", $true, $false, $false, $false, $false, $true, 1, $false, "64+18=82", 2) | Out-Null
$d.Content.Find.Execute("29-8=21", $true, $false, $false, $false, $false, $true, 1, $false, "8+45=53", 2) | Out-Null
$d.Content.Find.Execute("46+7=53", $true, $false, $false, $false, $false, $true, 1, $false, "37-1=36", 2) | Out-Null
$d.Content.Find.Execute("28+70=98", $true, $false, $false, $false, $false, $true, 1, $false, "4+3=7", 2) | Out-Null
$d.Content.Find.Execute("14+58=72", $true, $false, $false, $false, $false, $true, 1, $false, "12+66=78", 2) | Out-Null
$d.Content.Find.Execute("47-8=39", $true, $false, $false, $false, $false, $true, 1, $false, "26-19=7", 2) | Out-Null
$d.Content.Find.Execute("59-22=37", $true, $false, $false, $false, $false, $true, 1, $false, "43-11=32", 2) | Out-Null
$d.Content.Find.Execute("12+25=37", $true, $false, $false, $false, $false, $true, 1, $false, "48+3=51", 2) | Out-Null
$d.Content.Find.Execute("86-30=56", $true, $false, $false, $false, $false, $true, 1, $false, "16+9=25", 2) | Out-Null
$d.Content.Find.Execute("91-32=59", $true, $false, $false, $false, $false, $true, 1, $false, "96-76=20", 2) | Out-Null
$d.Content.Find.Execute("30+20=50", $true, $false, $false, $false, $false, $true, 1, $false, "58+31=89", 2) | Out-Null
$d.Content.Find.Execute("90-30=60", $true, $false, $false, $false, $false, $true, 1, $false, "7+27=34", 2) | Out-Null
$d.Content.Find.Execute("6+30=36", $true, $false, $false, $false, $false, $true, 1, $false, "75-19=56", 2) | Out-Null
$d.Content.Find.Execute("81+5=86", $true, $false, $false, $false, $false, $true, 1, $false, "29-23=6", 2) | Out-Null
$d.Content.Find.Execute("39-38=1", $true, $false, $false, $false, $false, $true, 1, $false, "40-20=20", 2) | Out-Null
$d.Content.Find.Execute("86-0=86", $true, $false, $false, $false, $false, $true, 1, $false, "49+40=89", 2) | Out-Null
$d.Content.Find.Execute("92-87=5", $true, $false, $false, $false, $false, $true, 1, $false, "55+35=90", 2) | Out-Null
$d.Content.Find.Execute("13-6=7", $true, $false, $false, $false, $false, $true, 1, $false, "12-7=5", 2) | Out-Null
$d.Content.Find.Execute("82+15=97", $true, $false, $false, $false, $false, $true, 1, $false, "89-72=17", 2) | Out-Null
$d.Content.Find.Execute("37+37=74", $true, $false, $false, $false, $false, $true, 1, $false, "23+37=60", 2) | Out-Null
$d.Content.Find.Execute("8+81=89", $true, $false, $false, $false, $false, $true, 1, $false, "56-18=38", 2) | Out-Null
$d.Content.Find.Execute("77-38=39", $true, $false, $false, $false, $false, $true, 1, $false, "78+0=78", 2) | Out-Null
$d.Content.Find.Execute("75-70=5", $true, $false, $false, $false, $false, $true, 1, $false, "58-27=31", 2) | Out-Null
$d.Content.Find.Execute("90-54=36", $true, $false, $false, $false, $false, $true, 1, $false, "70-1=69", 2) | Out-Null
$d.Content.Find.Execute("36+38=74", $true, $false, $false, $false, $false, $true, 1, $false, "26+21=47", 2) | Out-Null
$d.Content.Find.Execute("36+58=94", $true, $false, $false, $false, $false, $true, 1, $false, "7+24=31", 2) | Out-Null
$d.Content.Find.Execute("67-12=55", $true, $false, $false, $false, $false, $true, 1, $false, "11+6=17", 2) | Out-Null
$d.Content.Find.Execute("95-55=40", $true, $false, $false, $false, $false, $true, 1, $false, "83-19=64", 2) | Out-Null
$d.Content.Find.Execute("64+5=69", $true, $false, $false, $false, $false, $true, 1, $false, "66-28=38", 2) | Out-Null
$d.Content.Find.Execute("29+15=44", $true, $false, $false, $false, $false, $true, 1, $false, "16+53=69", 2) | Out-Null
$d.Content.Find.Execute("37-23=14", $true, $false, $false, $false, $false, $true, 1, $false, "67-61=6", 2) | Out-Null
$d.Content.Find.Execute("78+6=84", $true, $false, $false, $false, $false, $true, 1, $false, "39+58=97", 2) | Out-Null
$d.Content.Find.Execute("13+65=78", $true, $false, $false, $false, $false, $true, 1, $false, "54+35=89", 2) | Out-Null
$d.Content.Find.Execute("25+48=73", $true, $false, $false, $false, $false, $true, 1, $false, "73-38=35", 2) | Out-Null
$d.Content.Find.Execute("51-48=3", $true, $false, $false, $false, $false, $true, 1, $false, "91-8=83", 2) | Out-Null
$d.Content.Find.Execute("27-6=21", $true, $false, $false, $false, $false, $true, 1, $false, "30+44=74", 2) | Out-Null
$d.Content.Find.Execute("90-81=9", $true, $false, $false, $false, $false, $true, 1, $false, "20-5=15", 2) | Out-Null
$d.Content.Find.Execute("9+76=85", $true, $false, $false, $false, $false, $true, 1, $false, "85-52=33", 2) | Out-Null
$d.Content.Find.Execute("97-7=90", $true, $false, $false, $false, $false, $true, 1, $false, "95-11=84", 2) | Out-Null
$d.Content.Find.Execute("13+67=80", $true, $false, $false, $false, $false, $true, 1, $false, "35-0=35", 2) | Out-Null
$d.Content.Find.Execute("98-85=13", $true, $false, $false, $false, $false, $true, 1, $false, "27+25=52", 2) | Out-Null
$d.Content.Find.Execute("39+28=67", $true, $false, $false, $false, $false, $true, 1, $false, "18-15=3", 2) | Out-Null
$d.Content.Find.Execute("91+8=99", $true, $false, $false, $false, $false, $true, 1, $false, "10+44=54", 2) | Out-Null
$d.Content.Find.Execute("92-79=13", $true, $false, $false, $false, $false, $true, 1, $false, "21+41=62", 2) | Out-Null
$d.Content.Find.Execute("26+12=38", $true, $false, $false, $false, $false, $true, 1, $false, "44+13=57", 2) | Out-Null
$d.Content.Find.Execute("31+56=87", $true, $false, $false, $false, $false, $true, 1, $false, "26+42=68", 2) | Out-Null
$d.Content.Find.Execute("75-27=48", $true, $false, $false, $false, $false, $true, 1, $false, "18+57=75", 2) | Out-Null
$d.Content.Find.Execute("56+21=77", $true, $false, $false, $false, $false, $true, 1, $false, "20+71=91", 2) | Out-Null
$d.Content.Find.Execute("94-78=16", $true, $false, $false, $false, $false, $true, 1, $false, "43+6=49", 2) | Out-Null
$d.Content.Find.Execute("95-65=30", $true, $false, $false, $false, $false, $true, 1, $false, "27+7=34", 2) | Out-Null
$d.Content.Find.Execute("91-37=54", $true, $false, $false, $false, $false, $true, 1, $false, "63-45=18", 2) | Out-Null
$d.Content.Find.Execute("92-89=3", $true, $false, $false, $false, $false, $true, 1, $false, "11-0=11", 2) | Out-Null
$d.Content.Find.Execute("17-0=17", $true, $false, $false, $false, $false, $true, 1, $false, "10-9=1", 2) | Out-Null
$d.Content.Find.Execute("7+23=30", $true, $false, $false, $false, $false, $true, 1, $false, "56-30=26", 2) | Out-Null
$d.Content.Find.Execute("79-53=26", $true, $false, $false, $false, $false, $true, 1, $false, "58+33=91", 2) | Out-Null
$d.Content.Find.Execute("24+20=44", $true, $false, $false, $false, $false, $true, 1, $false, "49-39=10", 2) | Out-Null
$d.Content.Find.Execute("29+39=68", $true, $false, $false, $false, $false, $true, 1, $false, "4+25=29", 2) | Out-Null
$d.Content.Find.Execute("52+46=98", $true, $false, $false, $false, $false, $true, 1, $false, "68+12=80", 2) | Out-Null
$d.Content.Find.Execute("25+55=80", $true, $false, $false, $false, $false, $true, 1, $false, "17+20=37", 2) | Out-Null
$d.Content.Find.Execute("1+38=39", $true, $false, $false, $false, $false, $true, 1, $false, "14+78=92", 2) | Out-Null
$d.Content.Find.Execute("12+63=75", $true, $false, $false, $false, $false, $true, 1, $false, "48-45=3", 2) | Out-Null
$d.Content.Find.Execute("74-5=69", $true, $false, $false, $false, $false, $true, 1, $false, "77+19=96", 2) | Out-Null
$d.Content.Find.Execute("10+40=50", $true, $false, $false, $false, $false, $true, 1, $false, "70-25=45", 2) | Out-Null
$d.Content.Find.Execute("27+64=91", $true, $false, $false, $false, $false, $true, 1, $false, "73-69=4", 2) | Out-Null
$d.Content.Find.Execute("93-13=80", $true, $false, $false, $false, $false, $true, 1, $false, "67+30=97", 2) | Out-Null
$d.Content.Find.Execute("67+27=94", $true, $false, $false, $false, $false, $true, 1, $false, "71-6=65", 2) | Out-Null
$d.Content.Find.Execute("2+45=47", $true, $false, $false, $false, $false, $true, 1, $false, "17+79=96", 2) | Out-Null
$d.Content.Find.Execute("39+24=63", $true, $false, $false, $false, $false, $true, 1, $false, "87-39=48", 2) | Out-Null
$d.Content.Find.Execute("89-20=69", $true, $false, $false, $false, $false, $true, 1, $false, "15+12=27", 2) | Out-Null
$d.Content.Find.Execute("92-58=34", $true, $false, $false, $false, $false, $true, 1, $false, "6+49=55", 2) | Out-Null
$d.Content.Find.Execute("49+25=74", $true, $false, $false, $false, $false, $true, 1, $false, "70-52=18", 2) | Out-Null
$d.Content.Find.Execute("83-5=78", $true, $false, $false, $false, $false, $true, 1, $false, "27-16=11", 2) | Out-Null
$d.Content.Find.Execute("41+19=60", $true, $false, $false, $false, $false, $true, 1, $false, "64-48=16", 2) | Out-Null
$d.Content.Find.Execute("10+76=86", $true, $false, $false, $false, $false, $true, 1, $false, "31+58=89", 2) | Out-Null
$d.Content.Find.Execute("92-90=2", $true, $false, $false, $false, $false, $true, 1, $false, "43+28=71", 2) | Out-Null
$d.Content.Find.Execute("51-1=50", $true, $false, $false, $false, $false, $true, 1, $false, "7+66=73", 2) | Out-Null
$d.Content.Find.Execute("11+63=74", $true, $false, $false, $false, $false, $true, 1, $false, "83+13=96", 2) | Out-Null
$d.Content.Find.Execute("25-21=4", $true, $false, $false, $false, $false, $true, 1, $false, "97-65=32", 2) | Out-Null

Write-Output "Done: applied 100 replacements"
